# Atualização automática da planilha via Apps Script 003
#
# 1) Rename the sheet EstoqueSolicitacao -> Lotes (propagates to the
#    _xlnm._FilterDatabase defined name automatically).
# 2) Update the cached "imported" values (P/Q columns, fed from the
#    __xludf.DUMMYFUNCTION/IMPORTRANGE placeholder formulas) to the
#    refreshed numbers pulled in by the Apps Script sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Lotes"

$updates = @(
    @{ Cell = "P11";  New = 257.0 },
    @{ Cell = "P20";  New = 220.0 },
    @{ Cell = "P25";  New = 392.0 },
    @{ Cell = "P27";  New = 297.0 },
    @{ Cell = "P30";  New = 513.0 },
    @{ Cell = "P41";  New = 266.0 },
    @{ Cell = "P55";  New = 186.0 },
    @{ Cell = "P68";  New = 195.0 },
    @{ Cell = "P71";  New = 225.0 },
    @{ Cell = "P74";  New = 208.0 },
    @{ Cell = "P98";  New = 562.0 },
    @{ Cell = "P99";  New = 432.0 },
    @{ Cell = "P100"; New = 195.0 },
    @{ Cell = "P103"; New = 245.0 },
    @{ Cell = "P105"; New = 728.0 },
    @{ Cell = "Q105"; New = 2.0 },
    @{ Cell = "P107"; New = 598.0 },
    @{ Cell = "P111"; New = 77.0 },
    @{ Cell = "P112"; New = 77.0 },
    @{ Cell = "P127"; New = 579.0 },
    @{ Cell = "Q127"; New = 38.0 },
    @{ Cell = "P227"; New = 138.0 },
    @{ Cell = "P256"; New = 292.0 },
    @{ Cell = "P275"; New = 985.0 },
    @{ Cell = "Q275"; New = 94.0 },
    @{ Cell = "P301"; New = 542.0 },
    @{ Cell = "Q301"; New = 3.0 },
    @{ Cell = "P311"; New = 460.0 },
    @{ Cell = "P313"; New = 236.0 }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $oldFormula = $rng.Formula
    # Formula looks like:
    #   =IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),259.0)
    # Replace only the trailing numeric literal (the cached computed value)
    # with the refreshed number, keeping the DUMMYFUNCTION wrapper intact.
    $newValueText = [string]$u.New
    $replacement = $newValueText + ")"
    $newFormula = [System.Text.RegularExpressions.Regex]::Replace($oldFormula, '(-?[0-9]+(\.[0-9]+)?)\)$', $replacement)
    $rng.Formula = $newFormula
}
